# Add a "Price" column (N) with a value for every stock row (commit: "Add price for all stocks").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header ---
$ws.Range("N1").Value = "Price"

# --- Currency number format used by every populated Price cell ---
$currencyFmt = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# --- Values for N8:N50 (row 7 is handled separately below - different font) ---
$values = @(
    31.44, 30.98, 39.299999999999997, 40.49, 36.9, 38.06, 49.5, 36.01, 40.61,
    31.04, 27.17, 18.28, 20.46, 19.079999999999998, 8.89, 15.72, 12.06, 28.11,
    23.96, 22.92, 30.51, 33.31, 37.340000000000003, 46.13, 40.119999999999997,
    37.700000000000003, 45.74, 48.66, 48.47, 55.5, 57.45, 49.94, 52.13, 37.47,
    49.87, 44.44, 67.34, 66.790000000000006, 74.209999999999994, 99.38, 87.76,
    84.07, 90.35
)

# 1) Give every Price cell (N2:N50) the plain black Calibri font used across the column.
#    Doing this as one pass over the whole column first means the later, per-cell
#    NumberFormat / font-name tweak for N7 only has to change ONE extra property.
$ws.Range("N2:N50").Font.Color = 0

# 2) Rows 2-6 (header/blank rows in the source data) only get the font - no value, no $ format.
#    (left blank on purpose, matching the source rows that have no earnings data yet)

# 3) Rows 8-50 get the value plus the currency format.
$row = 8
foreach ($v in $values) {
    $cell = $ws.Cells.Item($row, 14)
    $cell.Value = $v
    $cell.NumberFormat = $currencyFmt
    $row = $row + 1
}

# 4) Row 7 gets its value, the currency format, and the "Aptos Narrow" font (as in the source).
$n7 = $ws.Range("N7")
$n7.Value = 27.59
$n7.NumberFormat = $currencyFmt
$n7.Font.Name = "Aptos Narrow"

# --- Selection, matching the saved selection in the edited workbook ---
$ws.Range("N2:N50").Select()
